# New cortisol analyses 02/03/2021
# Fill in Corti1-Corti4 (columns C-F) values for participants whose
# results had not yet been entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "9"   = @(0.04, 0.03, 0.02, 0.03)
    "30"  = @(0.33, 0.28, 0.1, 0.17)
    "58"  = @(0.17, 0.21)
    "61"  = @(0.12, 0.11, 0.11, 0.06)
    "62"  = @(0.1, 0.12, 0.11, 0.13)
    "63"  = @(0.26, 0.22, 0.24, 0.2)
    "72"  = @(0.19, 0.08, 0.06, 0.09)
    "73"  = @(0.02, 0.09, 0.08, 0.07000000000000001)
    "84"  = @(0.14, 0.15, 0.1, 0.08)
    "85"  = @(0.07000000000000001, 0.09, 0.39, 0.54)
    "86"  = @(0.42, 0.13, 0.08, 0.09)
    "87"  = @(0.34, 0.21, 0.34, 0.17)
    "88"  = @(0.17, 0.06, 0.08, 0.06)
    "89"  = @(0.09, 0.08, 0.05, 0.07000000000000001)
    "90"  = @(0.26, 0.17, 0.14, 0.12)
    "91"  = @(0.37, 0.23, 0.2, 0.25)
    "92"  = @(0.13, 0.15, 0.15, 0.1)
    "94"  = @(0.17, 0.17, 0.3, 0.31)
    "99"  = @(0.36, 0.49, 0.97, 0.87)
    "101" = @(0.26, 0.2, 0.26, 0.19)
    "102" = @(0.1, 0.13, 0.28, 0.39)
    "103" = @(0.25, 0.22, 0.22, 0.22)
    "105" = @(0.16, 0.15, 0.12, 0.09)
    "106" = @(0.23, 0.18, 0.15, 0.11)
    "112" = @(0.04, 0.08, 0.09, 0.06)
    "123" = @(0.17, 0.25, 0.21, 0.19)
    "125" = @(0.17, 0.17, 0.18, 0.14)
    "129" = @(0.05, 0.17, 0.02, 0.24)
    "130" = @(0.11, 0.08, 0.03, 0.03)
}

$columns = @("C", "D", "E", "F")

foreach ($rowNum in $updates.Keys) {
    $values = $updates[$rowNum]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $colLetter = $columns[$i]
        $ws.Range("$colLetter$rowNum").Value = $values[$i]
    }
}
